# Generate Report for Handoff
# Adds two newly-discovered handoff records (PNG dependency files tied to
# the existing markdown source) to the localization-status workbook:
#   - 594005b1-884a-40dc-9f3f-8e69d338b089.png  (dependency image, row 2)
#   - 5c90ffbf-4ace-41be-ac18-8282dba0c332.md   (new markdown source, row 3)
#   - a6857c3a-211d-4e70-809f-cbe6111e61c0.png  (dependency image, row 4)
# across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

function Set-LinkedCell {
    param($ws, $cellRef, $text, $url)
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text) | Out-Null
}

# ---------------------------------------------------------------------
# Source / target filenames used throughout the three sheets
# ---------------------------------------------------------------------
$mdFile   = "5c90ffbf-4ace-41be-ac18-8282dba0c332.md"
$png1     = "594005b1-884a-40dc-9f3f-8e69d338b089.png"
$png2     = "a6857c3a-211d-4e70-809f-cbe6111e61c0.png"

$zhTarget1 = "52dc8ca8da3401e11d54cf86e3c5f365f4b51ebb.png"
$zhTarget2 = "5c90ffbf-4ace-41be-ac18-8282dba0c332.ee980da7df0743a1719a5a4a6dd2e6c47d1cb463.zh-cn.xlf"
$zhTarget3 = "8c2e3080f6a27e7937cbc9b995c948dfc2f00f23.png"

$deTarget1 = "52dc8ca8da3401e11d54cf86e3c5f365f4b51ebb.png"
$deTarget2 = "5c90ffbf-4ace-41be-ac18-8282dba0c332.ee980da7df0743a1719a5a4a6dd2e6c47d1cb463.de-de.xlf"
$deTarget3 = "8c2e3080f6a27e7937cbc9b995c948dfc2f00f23.png"

$zhHandoffDateTime = "2016-03-20 04:50:20"
$deHandoffDateTime = "2016-03-20 04:50:23"
$overviewDate      = "2016-50-20 04:50:23"
$zeroDate          = "0001-01-01 00:00:00"

$statusReady   = "Ready for handoff"
$reasonDep     = "IsDependency"
$reasonInclude = "Include"
$dependsOnMd   = "e2e\$mdFile"

$srcBase = "https://github.com/OpenLocalizationTest/oltest/blob/84efe6f555b2c2bcd163a8514ddc93ae45f86e50/e2e"
$zhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/072277bce3336c67b00984feb88e1c09710b9732/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b839b096133e4d643ad474c00fecf00f45bfc04/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$png1Url = "$srcBase/$png1"
$mdUrl   = "$srcBase/$mdFile"
$png2Url = "$srcBase/$png2"

# =======================================================================
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

# NOTE: Range.Hyperlinks.Delete() on this engine clears every hyperlink on
# the worksheet (not just the addressed range), so we wipe once up front and
# re-add the complete, correctly-ordered set below.
$ov.Range("A1").Hyperlinks.Delete() | Out-Null

# Row 2 now describes the first PNG dependency (was the .md record before)
Set-LinkedCell $ov "A2" $png1 $png1Url
$ov.Range("B2").Value = $statusReady
$ov.Range("C2").Value = $statusReady
$ov.Range("D2").Value = $overviewDate

# Row 3: the markdown source file itself
Set-LinkedCell $ov "A3" $mdFile $mdUrl
$ov.Range("B3").Value = $statusReady
$ov.Range("C3").Value = $statusReady
$ov.Range("D3").Value = $overviewDate

# Row 4: the second PNG dependency
Set-LinkedCell $ov "A4" $png2 $png2Url
$ov.Range("B4").Value = $statusReady
$ov.Range("C4").Value = $statusReady
$ov.Range("D4").Value = $overviewDate

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A1").Hyperlinks.Delete() | Out-Null

# Row 2: PNG dependency #1
Set-LinkedCell $zh "A2" $png1 $png1Url
Set-LinkedCell $zh "B2" ".png" $png1Url
$zh.Range("C2").Value = $statusReady
Set-LinkedCell $zh "D2" $zhTarget1 "$zhBase/$zhTarget1"
$zh.Range("E2").Value = $zhHandoffDateTime
$zh.Range("H2").Value = $zeroDate
$zh.Range("I2").Value = $reasonDep
$zh.Range("J2").Value = $dependsOnMd

# Row 3: markdown source
Set-LinkedCell $zh "A3" $mdFile $mdUrl
Set-LinkedCell $zh "B3" ".md" $mdUrl
$zh.Range("C3").Value = $statusReady
Set-LinkedCell $zh "D3" $zhTarget2 "$zhBase/$zhTarget2"
$zh.Range("E3").Value = $zhHandoffDateTime
$zh.Range("H3").Value = $zeroDate
$zh.Range("I3").Value = $reasonInclude

# Row 4: PNG dependency #2
Set-LinkedCell $zh "A4" $png2 $png2Url
Set-LinkedCell $zh "B4" ".png" $png2Url
$zh.Range("C4").Value = $statusReady
Set-LinkedCell $zh "D4" $zhTarget3 "$zhBase/$zhTarget3"
$zh.Range("E4").Value = $zhHandoffDateTime
$zh.Range("H4").Value = $zeroDate
$zh.Range("I4").Value = $reasonDep
$zh.Range("J4").Value = $dependsOnMd

# =======================================================================
# Sheet "de-de"
# =======================================================================
$de = $wb.Worksheets.Item("de-de")
$de.Range("A1").Hyperlinks.Delete() | Out-Null

# Row 2: PNG dependency #1
Set-LinkedCell $de "A2" $png1 $png1Url
Set-LinkedCell $de "B2" ".png" $png1Url
$de.Range("C2").Value = $statusReady
Set-LinkedCell $de "D2" $deTarget1 "$deBase/$deTarget1"
$de.Range("E2").Value = $deHandoffDateTime
$de.Range("H2").Value = $zeroDate
$de.Range("I2").Value = $reasonDep
$de.Range("J2").Value = $dependsOnMd

# Row 3: markdown source
Set-LinkedCell $de "A3" $mdFile $mdUrl
Set-LinkedCell $de "B3" ".md" $mdUrl
$de.Range("C3").Value = $statusReady
Set-LinkedCell $de "D3" $deTarget2 "$deBase/$deTarget2"
$de.Range("E3").Value = $deHandoffDateTime
$de.Range("H3").Value = $zeroDate
$de.Range("I3").Value = $reasonInclude

# Row 4: PNG dependency #2
Set-LinkedCell $de "A4" $png2 $png2Url
Set-LinkedCell $de "B4" ".png" $png2Url
$de.Range("C4").Value = $statusReady
Set-LinkedCell $de "D4" $deTarget3 "$deBase/$deTarget3"
$de.Range("E4").Value = $deHandoffDateTime
$de.Range("H4").Value = $zeroDate
$de.Range("I4").Value = $reasonDep
$de.Range("J4").Value = $dependsOnMd
